# adelanto del informe: +http/https +firewall
#
# Applies the documented edit:
#   1. "Los cuatro primeros servicios" -> "Los tres primeros servicios"
#      (and relocates the internal "_GoBack" bookmark to sit right after
#      "Los tres " as in the target OOXML).
#   2. Appends two new sentences (HTTP/HTTPS explanation) to the
#      "Los servidores web " paragraph.
#   3. Adds one extra blank paragraph after that paragraph.
#   4. Replaces the "Firewall -> nivel aplicación(7mo)" paragraph with a
#      new paragraph describing the Secretaria de Hacienda firewall.
#   5. Deletes the "Cables -> físico" paragraph entirely.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "cuatro" -> "tres", then move the _GoBack bookmark so it sits right
#    before "primeros servicios" (i.e. right after "Los tres ").
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "Los cuatro primeros servicios"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "Los tres primeros servicios"
$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null

$implParaRange = $d.Paragraphs.Item(58).Range.Duplicate
$bmFind = $implParaRange.Find
$bmFind.ClearFormatting()
$bmFind.Text = "Los tres "
$bmFind.Execute() | Out-Null
$implParaRange.Collapse(0)

$d.Bookmarks.Add("_GoBack", $implParaRange)

# ---------------------------------------------------------------------
# 2) Extend the "Los servidores web " paragraph with the HTTP/HTTPS text.
# ---------------------------------------------------------------------
$webPara = $d.Paragraphs.Item(59)
$webRange = $webPara.Range
$webRange.End = $webRange.End - 1
$webRange.Collapse(0)
$webRange.InsertAfter("están divididos en dos grupos: los que trabajan con protocolo HTTP y los que lo hacen con protocolo HTTPS, ambos protocolos de capa de aplicación. A través de estos protocolos, los clientes que lo requieran pueden obtener los conjuntos de archivos que hacen a una página web: HTML, CSS, JS, etc. ")
$webRange.Collapse(0)
$webRange.InsertAfter("La principal diferencia entre HTTPS y HTTP es que en el primero la información que va del cliente al servidor viaja cifrada, lo que permite que frente a una eventual captura de estos datos no sea posible obtener directamente la información que contienen, en cambio en el segundo los datos viajan en texto plano.")

# ---------------------------------------------------------------------
# 3) One extra empty paragraph right after it.
# ---------------------------------------------------------------------
$d.Paragraphs.Item(59).Range.InsertParagraphAfter()

# ---------------------------------------------------------------------
# 4) Replace the "Firewall -> nivel aplicación(7mo)" paragraph.
#    (Locate the paragraph object itself rather than relying on a
#    Find-match sub-range, whose .Paragraphs collection does not expand
#    back out to the full enclosing paragraph.)
# ---------------------------------------------------------------------
$fwPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Firewall -> nivel")) {
        $fwPara = $p
        break
    }
}
$fwParaRange = $fwPara.Range
$fwParaRange.End = $fwParaRange.End - 1
$fwParaRange.Text = "En el servidor web de la Secretaria de hacienda, además de utilizar el protocolo HTTPS, tiene configurado un firewall que permite conectarse a ese servidor solo a clientes que tengan determinadas IPs especificadas en la configuración. Como tiene que filtrar paquetes IP, este firewall trabaja en capa 3."

# ---------------------------------------------------------------------
# 5) Delete the "Cables -> físico" paragraph entirely.
# ---------------------------------------------------------------------
$cablesPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Cables ->")) {
        $cablesPara = $p
        break
    }
}
$cablesPara.Range.Delete()
